$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 17
$ws.Range('F7').Value = 898
$ws.Range('F8').Value = 558
$ws.Range('F9').Value = 2396
$ws.Range('F10').Value = 689
$ws.Range('F12').Value = 556
$ws.Range('F14').Value = 318
$ws.Range('F15').Value = 194
$ws.Range('F16').Value = 512
$ws.Range('F17').Value = 2098
$ws.Range('F19').Value = 691
$ws.Range('F21').Value = 2598
$ws.Range('F27').Value = 1743
$ws.Range('F30').Value = 511
$ws.Range('F33').Value = 4520
$ws.Range('F34').Value = 80

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F2').Value = 404
$ws.Range('F3').Value = 381
$ws.Range('F4').Value = 8
$ws.Range('F14').Value = 314
$ws.Range('F17').Value = 151
$ws.Range('F24').Value = 216
$ws.Range('F26').Value = 241

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 1406
$ws.Range('F7').Value = 158

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 1406
$ws.Range('F6').Value = 381
$ws.Range('F7').Value = 17
$ws.Range('F15').Value = 898
$ws.Range('F16').Value = 558
$ws.Range('F17').Value = 2396
$ws.Range('F18').Value = 689
$ws.Range('F20').Value = 556
$ws.Range('F22').Value = 318
$ws.Range('F24').Value = 194
$ws.Range('F25').Value = 314
$ws.Range('F26').Value = 512
$ws.Range('F27').Value = 2098
$ws.Range('F29').Value = 691
$ws.Range('B30').Value = '''2024-04-20'
$ws.Range('C30').Value = '上海· 茅原实里动漫交响音乐会'
$ws.Range('D30').Value = '东大名路889号 友邦大剧院'
$ws.Range('E30').Value = '2024.04.20 19:30-04.20 21:00'
$ws.Range('F30').Value = 151
$ws.Range('G30').Value = 380
$ws.Range('H30').Value = 'https://show.bilibili.com/platform/detail.html?id=81703'
$ws.Range('I30').Value = '//i2.hdslb.com/bfs/openplatform/202402/yiVaqJVK1707016321221.jpeg'
$ws.Range('C31').Value = '上海·K-9AL动漫展'
$ws.Range('D31').Value = '市真南路1199弄1号 智创TOP综合体产城'
$ws.Range('E31').Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range('F31').Value = 0
$ws.Range('G31').Value = 78
$ws.Range('H31').Value = 'https://show.bilibili.com/platform/detail.html?id=83376'
$ws.Range('I31').Value = '//i2.hdslb.com/bfs/openplatform/202403/zc80cfXW1711374771275.png'
$ws.Range('C32').Value = '上海·Virtual Shanghai Anime Exhibition魔都虚幻世界二次元1.0'
$ws.Range('D32').Value = '中山北路3300号环球港购物中心4楼 上海世嘉都市乐园'
$ws.Range('F32').Value = 2598
$ws.Range('G32').Value = 69
$ws.Range('H32').Value = 'https://show.bilibili.com/platform/detail.html?id=81865'
$ws.Range('I32').Value = '//i0.hdslb.com/bfs/openplatform/202403/C8G3AOLM1709870001354.jpeg'
$ws.Range('C33').Value = '上海·冰兔2024线下live《过去和未来》'
$ws.Range('D33').Value = '重庆南路308号3楼 上海MAO LIVEHOUSE'
$ws.Range('E33').Value = '2024.04.20 13:00-04.20 15:00'
$ws.Range('F33').Value = 269
$ws.Range('G33').Value = 198
$ws.Range('H33').Value = 'https://show.bilibili.com/platform/detail.html?id=81654'
$ws.Range('I33').Value = '//i1.hdslb.com/bfs/openplatform/202402/OEHnMZmi1706851347869.jpeg'
$ws.Range('C34').Value = '上海·心苑主题现场·《Husky Go × 阿君归来》联动主题签名会'
$ws.Range('D34').Value = '长宁路890号玫瑰坊B2-16号 Husky Go 哈士奇体验馆'
$ws.Range('E34').Value = '2024.04.20 13:50-04.20 18:00'
$ws.Range('F34').Value = 35
$ws.Range('G34').Value = 106
$ws.Range('H34').Value = 'https://show.bilibili.com/platform/detail.html?id=83114'
$ws.Range('I34').Value = '//i2.hdslb.com/bfs/openplatform/202403/Hlb7aPAX1710989888372.jpeg'
$ws.Range('C35').Value = '上海·环绕银河动漫游戏嘉年华-音你跃动'
$ws.Range('D35').Value = '逸仙路1328弄 新业坊'
$ws.Range('E35').Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range('F35').Value = 22
$ws.Range('G35').Value = 68
$ws.Range('H35').Value = 'https://show.bilibili.com/platform/detail.html?id=82839'
$ws.Range('I35').Value = '//i1.hdslb.com/bfs/openplatform/202403/lui5Ed5W1710385702545.png'
$ws.Range('C36').Value = '上海·第五十六届妖漫动漫游戏展'
$ws.Range('D36').Value = '共和新路2188号 上海久光中心'
$ws.Range('E36').Value = '2024.04.20 10:00-04.20 17:00'
$ws.Range('F36').Value = 510
$ws.Range('G36').Value = 77.7
$ws.Range('H36').Value = 'https://show.bilibili.com/platform/detail.html?id=83298'
$ws.Range('I36').Value = '//i0.hdslb.com/bfs/openplatform/202403/PDYSzPVC1711255759583.jpeg'
$ws.Range('B37').Value = '''2024-04-21'
$ws.Range('C37').Value = '上海·今泉爱夏  巡演'
$ws.Range('D37').Value = '瑞虹路188号3楼 Modernsky Lab'
$ws.Range('E37').Value = '2024.04.21 20:00-04.21 21:30'
$ws.Range('F37').Value = 58
$ws.Range('G37').Value = 328
$ws.Range('H37').Value = 'https://show.bilibili.com/platform/detail.html?id=81891'
$ws.Range('I37').Value = '//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg'
$ws.Range('B38').Value = '''2024-04-24'
$ws.Range('C38').Value = '上海·「NIJISANJI EN x animate cafe」'
$ws.Range('D38').Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws.Range('E38').Value = '2024.04.24 00:00-05.22 23:59'
$ws.Range('F38').Value = 158
$ws.Range('G38').Value = 30
$ws.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=83223'
$ws.Range('I38').Value = '//i0.hdslb.com/bfs/openplatform/202403/LzJJK9lc1711096202393.jpeg'
$ws.Range('B39').Value = '''2024-04-26'
$ws.Range('C39').Value = '上海·「angela LIVE 2024」in SHANGHAI'
$ws.Range('D39').Value = '万航渡后路19号3楼 瓦肆VAS SHANGHAI'
$ws.Range('E39').Value = '2024.04.26 19:00-04.26 20:30'
$ws.Range('F39').Value = 1758
$ws.Range('G39').Value = 480
$ws.Range('H39').Value = 'https://show.bilibili.com/platform/detail.html?id=82039'
$ws.Range('I39').Value = '//i2.hdslb.com/bfs/openplatform/202402/H9L22d9R1708678603570.jpeg'
$ws.Range('B40').Value = '''2024-04-27'
$ws.Range('C40').Value = '上海·坏孩纸物语第39届动漫节'
$ws.Range('D40').Value = '曹杨路2033号 普陀绿地缤纷城'
$ws.Range('E40').Value = '2024.04.27 10:00-04.27 17:00'
$ws.Range('F40').Value = 507
$ws.Range('G40').Value = 55.5
$ws.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=83300'
$ws.Range('I40').Value = '//i1.hdslb.com/bfs/openplatform/202403/WqMFX7w01711261080799.png'
$ws.Range('B41').Value = '''2024-05-01'
$ws.Range('C41').Value = '上海·第三届Redamancy动漫游戏嘉年华'
$ws.Range('D41').Value = '中山北路3300号4楼L4001号 环球港上海世嘉都市乐园'
$ws.Range('E41').Value = '2024.05.01 10:00-05.03 17:00'
$ws.Range('F41').Value = 1743
$ws.Range('G41').Value = 60
$ws.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=82017'
$ws.Range('I41').Value = '//i1.hdslb.com/bfs/openplatform/202402/UJkFbmo91708657659067.png'
$ws.Range('C42').Value = '上海·街舞音乐剧《时光代理人：法则游戏》'
$ws.Range('D42').Value = '牛庄路704号 中国大戏院'
$ws.Range('E42').Value = '2024.05.01 19:30-05.19 21:00'
$ws.Range('F42').Value = 216
$ws.Range('G42').Value = 188
$ws.Range('H42').Value = 'https://show.bilibili.com/platform/detail.html?id=82995'
$ws.Range('I42').Value = '//i1.hdslb.com/bfs/openplatform/202403/p9ZC2azX1710816437198.png'
$ws.Range('F44').Value = 511
$ws.Range('F47').Value = 4520
$ws.Range('F48').Value = 80
